$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(
    @{ Row = 208; Time = "2023-12-11 20:28:25"; Cost = 0.0004 },
    @{ Row = 209; Time = "2023-12-11 20:29:25"; Cost = 0.003600000000000001 },
    @{ Row = 210; Time = "2023-12-11 20:29:47"; Cost = 0.0016 },
    @{ Row = 211; Time = "2023-12-11 20:29:55"; Cost = 0.0004 },
    @{ Row = 212; Time = "2023-12-11 20:30:05"; Cost = 0.0004 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Time
    $ws.Cells.Item($r.Row, 2).Value = $r.Cost
}
